$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: shift the quarterly summary rows down by one and insert
#    the new 2022-Q4 row at the top (row 2), pushing 2021-Q3 into a new row 7.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Give the new bottom row (A7) the same style as the existing column-A cells
# before writing values into it.
$totals.Range("A7").Value = 5
$totals.Range("A6").Copy()
$totals.Range("A7").PasteSpecial(-4122)

$totalsData = @(
    @(0, "2022-Q4", 14, 0.69),
    @(1, "2022-Q3", 12, 0.38),
    @(2, "2022-Q2", 15, 0.69),
    @(3, "2022-Q1", 9, 1.41),
    @(4, "2021-Q4", 19, 2.69),
    @(5, "2021-Q3", 15, 4.87)
)

$r = 2
foreach ($row in $totalsData) {
    $totals.Range("A$r").Value = $row[0]
    $totals.Range("B$r").Value = $row[1]
    $totals.Range("C$r").Value = $row[2]
    $totals.Range("D$r").Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" worksheet right after "总计" (i.e. before
#    the worksheet that is currently "2022-Q3"), matching the formatting of
#    the existing quarterly fund-holding sheets.
# ---------------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($oldQ3)
$newSheet.Name = "2022-Q4"

# Copy the header row + column-A formatting from one of the existing
# quarterly sheets (use 2021-Q4, which has 19 data rows - plenty to cover
# the 14 rows we are about to write - so every destination cell has a
# same-column source cell to copy formatting from).
$q3Ref = $wb.Worksheets.Item(3)
$q3Ref.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$formatRef = $wb.Worksheets.Item(6)
$formatRef.Range("A2:A15").Copy()
$newSheet.Range("A2:A15").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundData = @(
    @(0, "002620", "中邮未来新蓝筹灵活配置混合", "11.76", "86.82", "3.40", "0.3998", 6),
    @(1, "015729", "朱雀碳中和三年持有期混合", "3.60", "49.87", "2.26", "0.0814", 9),
    @(2, "014175", "工银价值成长混合A", "2.13", "73.24", "2.97", "0.0633", 4),
    @(3, "519618", "银河君信灵活配置混合I", "2.93", "21.62", "1.09", "0.0319", 4),
    @(4, "012132", "华泰保兴价值成长混合A", "0.80", "81.95", "3.73", "0.0298", 10),
    @(5, "012430", "农银汇理瑞康6个月持有期混合", "1.16", "28.06", "1.38", "0.0160", 3),
    @(6, "519656", "银河灵活配置混合 - A", "0.41", "72.68", "3.75", "0.0154", 5),
    @(7, "519657", "银河灵活配置混合 - C", "0.34", "72.68", "3.75", "0.0128", 5),
    @(8, "008124", "中邮中证500指数增强C", "0.90", "93.17", "1.37", "0.0123", 9),
    @(9, "014176", "工银价值成长混合C", "0.37", "73.24", "2.97", "0.0110", 4),
    @(10, "519617", "银河君信灵活配置混合C", "0.38", "21.62", "1.09", "0.0041", 4),
    @(11, "519616", "银河君信灵活配置混合A", "0.34", "21.62", "1.09", "0.0037", 4),
    @(12, "590007", "中邮中证500指数增强A", "0.26", "93.17", "1.37", "0.0036", 9),
    @(13, "012177", "华泰保兴价值成长混合C", "0.08", "81.95", "3.73", "0.0030", 10)
)

$r = 2
foreach ($row in $fundData) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").NumberFormat = "@"
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").NumberFormat = "@"
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").NumberFormat = "@"
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("F$r").NumberFormat = "@"
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("G$r").NumberFormat = "@"
    $newSheet.Range("G$r").Value = $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}
